$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the edited columns as Text so that numeric-looking values
# ("261.34", "1.06%", "19", ...) are preserved verbatim as text, matching
# the source data which stores every cell as plain text.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "261.34"
$ws.Range("E2").Value = "1.06%"
$ws.Range("G2").Value = "19"
$ws.Range("D3").Value = "27.08"
$ws.Range("E3").Value = "0.75%"
$ws.Range("G3").Value = "19"
$ws.Range("D4").Value = "4.710"
$ws.Range("E4").Value = "1.03%"
$ws.Range("G4").Value = "19"
$ws.Range("D5").Value = "0.06190"
$ws.Range("E5").Value = "3.21%"
$ws.Range("G5").Value = "19"
$ws.Range("D6").Value = "6.712"
$ws.Range("E6").Value = "0.73%"
$ws.Range("G6").Value = "19"
$ws.Range("E7").Value = "-0.95%"
$ws.Range("G7").Value = "19"
$ws.Range("D8").Value = "0.9155"
$ws.Range("E8").Value = "-1.15%"
$ws.Range("G8").Value = "19"
$ws.Range("D9").Value = "0.1411"
$ws.Range("E9").Value = "1.56%"
$ws.Range("G9").Value = "19"
$ws.Range("D10").Value = "0.04606"
$ws.Range("E10").Value = "-2.43%"
$ws.Range("G10").Value = "19"
$ws.Range("E11").Value = "1.01%"
$ws.Range("G11").Value = "19"
$ws.Range("D12").Value = "0.03133"
$ws.Range("E12").Value = "1.33%"
$ws.Range("G12").Value = "19"
$ws.Range("D13").Value = "0.09049"
$ws.Range("E13").Value = "-0.89%"
$ws.Range("G13").Value = "19"
$ws.Range("D14").Value = "0.001538"
$ws.Range("E14").Value = "0.79%"
$ws.Range("G14").Value = "19"
$ws.Range("D15").Value = "0.0006158"
$ws.Range("E15").Value = "1.82%"
$ws.Range("G15").Value = "19"
$ws.Range("D16").Value = "0.006101"
$ws.Range("E16").Value = "-1.05%"
$ws.Range("G16").Value = "19"
$ws.Range("D17").Value = "3.457"
$ws.Range("E17").Value = "-0.02%"
$ws.Range("G17").Value = "19"
$ws.Range("D18").Value = "3.166"
$ws.Range("E18").Value = "0.36%"
$ws.Range("G18").Value = "19"
$ws.Range("E19").Value = "0.40%"
$ws.Range("G19").Value = "19"
$ws.Range("E20").Value = "-0.97%"
$ws.Range("G20").Value = "19"
$ws.Range("D21").Value = "0.1310"
$ws.Range("E21").Value = "1.71%"
$ws.Range("G21").Value = "19"
$ws.Range("D22").Value = "4.126"
$ws.Range("E22").Value = "-0.40%"
$ws.Range("G22").Value = "19"
$ws.Range("D23").Value = "0.04232"
$ws.Range("E23").Value = "-0.13%"
$ws.Range("G23").Value = "19"
$ws.Range("E24").Value = "-0.05%"
$ws.Range("G24").Value = "19"
$ws.Range("E25").Value = "-5.79%"
$ws.Range("G25").Value = "19"
$ws.Range("G26").Value = "19"
$ws.Range("G27").Value = "19"
$ws.Range("G28").Value = "19"
$ws.Range("G29").Value = "19"
$ws.Range("G30").Value = "19"
$ws.Range("G31").Value = "19"
$ws.Range("G32").Value = "19"
$ws.Range("G33").Value = "19"
$ws.Range("G34").Value = "19"
$ws.Range("G35").Value = "19"
$ws.Range("G36").Value = "19"
$ws.Range("G37").Value = "19"
$ws.Range("G38").Value = "19"
$ws.Range("G39").Value = "19"
$ws.Range("D40").Value = "0.03994"
$ws.Range("E40").Value = "4.11%"
$ws.Range("G40").Value = "19"
$ws.Range("D41").Value = "0.1113"
$ws.Range("E41").Value = "-0.04%"
$ws.Range("G41").Value = "19"
$ws.Range("E42").Value = "7.47%"
$ws.Range("G42").Value = "19"
$ws.Range("G43").Value = "19"
$ws.Range("E44").Value = "-8.17%"
$ws.Range("G44").Value = "19"
$ws.Range("D45").Value = "0.00005161"
$ws.Range("E45").Value = "1.30%"
$ws.Range("G45").Value = "19"
$ws.Range("E46").Value = "0.07%"
$ws.Range("G46").Value = "19"
$ws.Range("G47").Value = "19"
$ws.Range("G48").Value = "19"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "0.07%"
$ws.Range("G49").Value = "19"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "0.07%"
$ws.Range("G50").Value = "19"
$ws.Range("G51").Value = "19"

# Restore the default (unstyled) look now that the text values are locked in.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
